$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.893.43'
$ws.Range("E2").Value = '  +2.47%  '

$ws.Range("D3").Value = '2.619.89'
$ws.Range("E3").Value = '  +0.87%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.38%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.584'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.37%  '

$ws.Range("E9").Value = '  +7.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.406'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.48%  '

$ws.Range("E11").Value = '  -0.63%  '

$ws.Range("E12").Value = '  +2.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.03'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000188'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +22.61%  '

$ws.Range("D15").Value = '3.088.36'
$ws.Range("E15").Value = '  +0.83%  '

$ws.Range("D16").Value = '64.840.25'
$ws.Range("E16").Value = '  +2.67%  '

$ws.Range("D17").Value = '2.620.11'
$ws.Range("E17").Value = '  +1.77%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '356.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.34%  '

$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.43'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.97%  '

$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.66'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.06%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.55%  '

$ws.Range("E27").Value = '  +2.50%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.39%  '

$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("D30").Value = '0.0₃0942'
$ws.Range("E30").Value = '  +12.56%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.47%  '

$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '524.06'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.64%  '

$ws.Range("E33").Value = '  +2.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.14%  '

$ws.Range("E35").Value = '  +5.19%  '

$ws.Range("E36").Value = '  +3.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '162.44'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.55%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '41.99'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '164.06'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.36%  '

$ws.Range("E44").Value = '  +4.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0611'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.648'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0260'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.77%  '

$ws.Range("E50").Value = '  +2.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.39'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.35%  '
